# Scheduled runner update: refresh market-board pricing snapshots
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Cells.Item(74, 8).Value = 6328.9165
$ws.Cells.Item(74, 9).Value = 5788.778
$ws.Cells.Item(74, 10).Value = 7949.3335
$ws.Cells.Item(74, 11).Value = 5788.778
$ws.Cells.Item(74, 12).Value = 7949.3335
$ws.Cells.Item(74, 13).Value = -4852.778
$ws.Cells.Item(74, 14).Value = -9821.333500000001
# Row 77
$ws.Cells.Item(77, 8).Value = 6328.9165
$ws.Cells.Item(77, 9).Value = 5788.778
$ws.Cells.Item(77, 10).Value = 7949.3335
$ws.Cells.Item(77, 11).Value = 28943.89
$ws.Cells.Item(77, 12).Value = 39746.6675
$ws.Cells.Item(77, 13).Value = -24263.89
$ws.Cells.Item(77, 14).Value = -49106.6675
# Row 98
$ws.Cells.Item(98, 8).Value = 7941
$ws.Cells.Item(98, 9).Value = 7958.5
$ws.Cells.Item(98, 10).Value = 7929.3335
$ws.Cells.Item(98, 11).Value = 7958.5
$ws.Cells.Item(98, 12).Value = 7929.3335
$ws.Cells.Item(98, 13).Value = -6460.5
$ws.Cells.Item(98, 14).Value = -10925.3335
# Row 112
$ws.Cells.Item(112, 8).Value = 1211805.4
$ws.Cells.Item(112, 10).Value = 1557093.2
$ws.Cells.Item(112, 12).Value = 4671279.6
$ws.Cells.Item(112, 14).Value = -4673495.6
# Row 122
$ws.Cells.Item(122, 8).Value = 7941
$ws.Cells.Item(122, 9).Value = 7958.5
$ws.Cells.Item(122, 10).Value = 7929.3335
$ws.Cells.Item(122, 11).Value = 23875.5
$ws.Cells.Item(122, 12).Value = 23788.0005
$ws.Cells.Item(122, 13).Value = -21425.5
$ws.Cells.Item(122, 14).Value = -28688.0005
# Row 129
$ws.Cells.Item(129, 8).Value = 2370.74
$ws.Cells.Item(129, 9).Value = 1023.125
$ws.Cells.Item(129, 10).Value = 2487.9238
$ws.Cells.Item(129, 11).Value = 3069.375
$ws.Cells.Item(129, 12).Value = 7463.7714
$ws.Cells.Item(129, 13).Value = 1930.625
$ws.Cells.Item(129, 14).Value = -17463.7714
# Row 133
$ws.Cells.Item(133, 8).Value = 89999
$ws.Cells.Item(133, 10).Value = 89999
$ws.Cells.Item(133, 12).Value = 89999
$ws.Cells.Item(133, 14).Value = -100119
# Row 134
$ws.Cells.Item(134, 8).Value = 50000
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
# Row 137
$ws.Cells.Item(137, 8).Value = 13025.454
$ws.Cells.Item(137, 9).Value = 4571.6665
$ws.Cells.Item(137, 11).Value = 13714.9995
$ws.Cells.Item(137, 13).Value = -11164.9995
# Row 138
$ws.Cells.Item(138, 8).Value = 6644.077
$ws.Cells.Item(138, 9).Value = 4000
$ws.Cells.Item(138, 10).Value = 6864.4165
$ws.Cells.Item(138, 11).Value = 12000
$ws.Cells.Item(138, 12).Value = 20593.2495
$ws.Cells.Item(138, 13).Value = -6860
$ws.Cells.Item(138, 14).Value = -30873.2495

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 34568304
$ws.Cells.Item(32, 9).Value = 37052996
$ws.Cells.Item(32, 11).Value = 37052996
$ws.Cells.Item(32, 13).Value = -37052709
# Row 51
$ws.Cells.Item(51, 8).Value = 44999
$ws.Cells.Item(51, 10).Value = 44999
$ws.Cells.Item(51, 12).Value = 44999
$ws.Cells.Item(51, 14).Value = -46511
# Row 63
$ws.Cells.Item(63, 8).Value = 27431.143
$ws.Cells.Item(63, 9).Value = 27431.143
$ws.Cells.Item(63, 11).Value = 27431.143
$ws.Cells.Item(63, 13).Value = -26745.143
# Row 66
$ws.Cells.Item(66, 8).Value = 27431.143
$ws.Cells.Item(66, 9).Value = 27431.143
$ws.Cells.Item(66, 11).Value = 137155.715
$ws.Cells.Item(66, 13).Value = -133723.715
# Row 74
$ws.Cells.Item(74, 8).Value = 2095.3713
$ws.Cells.Item(74, 9).Value = 980.9091
$ws.Cells.Item(74, 10).Value = 2606.1667
$ws.Cells.Item(74, 11).Value = 980.9091
$ws.Cells.Item(74, 12).Value = 2606.1667
$ws.Cells.Item(74, 13).Value = -106.9091
$ws.Cells.Item(74, 14).Value = -4354.1667
# Row 77
$ws.Cells.Item(77, 8).Value = 2095.3713
$ws.Cells.Item(77, 9).Value = 980.9091
$ws.Cells.Item(77, 10).Value = 2606.1667
$ws.Cells.Item(77, 11).Value = 4904.5455
$ws.Cells.Item(77, 12).Value = 13030.8335
$ws.Cells.Item(77, 13).Value = -536.5455000000002
$ws.Cells.Item(77, 14).Value = -21766.8335
# Row 132
$ws.Cells.Item(132, 8).Value = 50002544
$ws.Cells.Item(132, 9).Value = 2878
$ws.Cells.Item(132, 10).Value = 83335656
$ws.Cells.Item(132, 11).Value = 8634
$ws.Cells.Item(132, 12).Value = 250006968
$ws.Cells.Item(132, 13).Value = -6104
$ws.Cells.Item(132, 14).Value = -250012028

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 4918.4546
$ws.Cells.Item(86, 9).Value = 5340.6
$ws.Cells.Item(86, 10).Value = 4566.6665
$ws.Cells.Item(86, 11).Value = 5340.6
$ws.Cells.Item(86, 12).Value = 4566.6665
$ws.Cells.Item(86, 13).Value = -4217.6
$ws.Cells.Item(86, 14).Value = -6812.6665
# Row 89
$ws.Cells.Item(89, 8).Value = 4918.4546
$ws.Cells.Item(89, 9).Value = 5340.6
$ws.Cells.Item(89, 10).Value = 4566.6665
$ws.Cells.Item(89, 11).Value = 26703
$ws.Cells.Item(89, 12).Value = 22833.3325
$ws.Cells.Item(89, 13).Value = -21087
$ws.Cells.Item(89, 14).Value = -34065.3325
# Row 134
$ws.Cells.Item(134, 8).Value = 53040636
$ws.Cells.Item(134, 9).Value = 31260582
$ws.Cells.Item(134, 11).Value = 93781746
$ws.Cells.Item(134, 13).Value = -93779211

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5684318
$ws.Cells.Item(31, 9).Value = 6062805.5
$ws.Cells.Item(31, 10).Value = 7000
$ws.Cells.Item(31, 11).Value = 6062805.5
$ws.Cells.Item(31, 12).Value = 7000
$ws.Cells.Item(31, 13).Value = -6062510.5
$ws.Cells.Item(31, 14).Value = -7590
# Row 34
$ws.Cells.Item(34, 8).Value = 5684318
$ws.Cells.Item(34, 9).Value = 6062805.5
$ws.Cells.Item(34, 10).Value = 7000
$ws.Cells.Item(34, 11).Value = 6062805.5
$ws.Cells.Item(34, 12).Value = 7000
$ws.Cells.Item(34, 13).Value = -6062603.5
$ws.Cells.Item(34, 14).Value = -7404
# Row 58
$ws.Cells.Item(58, 8).Value = 2313.6667
$ws.Cells.Item(58, 9).Value = 2396.4
$ws.Cells.Item(58, 10).Value = 1900
$ws.Cells.Item(58, 11).Value = 2396.4
$ws.Cells.Item(58, 12).Value = 1900
$ws.Cells.Item(58, 13).Value = -2193.4
$ws.Cells.Item(58, 14).Value = -2306
# Row 134
$ws.Cells.Item(134, 8).Value = 2858662
$ws.Cells.Item(134, 9).Value = 1547
$ws.Cells.Item(134, 10).Value = 16668051
$ws.Cells.Item(134, 11).Value = 4641
$ws.Cells.Item(134, 12).Value = 50004153
$ws.Cells.Item(134, 13).Value = -2106
$ws.Cells.Item(134, 14).Value = -50009223
# Row 136
$ws.Cells.Item(136, 8).Value = 2313.6667
$ws.Cells.Item(136, 9).Value = 2396.4
$ws.Cells.Item(136, 10).Value = 1900
$ws.Cells.Item(136, 11).Value = 7189.200000000001
$ws.Cells.Item(136, 12).Value = 5700
$ws.Cells.Item(136, 13).Value = -4639.200000000001
$ws.Cells.Item(136, 14).Value = -10800
# Row 141
$ws.Cells.Item(141, 8).Value = 333406.88
$ws.Cells.Item(141, 10).Value = 395175.84
$ws.Cells.Item(141, 12).Value = 395175.84
$ws.Cells.Item(141, 14).Value = -405535.84

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 739727.9399999999
$ws.Cells.Item(131, 10).Value = 828486.75
$ws.Cells.Item(131, 12).Value = 2485460.25
$ws.Cells.Item(131, 14).Value = -2495540.25

$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Cells.Item(47, 8).Value = 30999.5
$ws.Cells.Item(47, 10).Value = 30999.5
$ws.Cells.Item(47, 12).Value = 30999.5
$ws.Cells.Item(47, 14).Value = -32135.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 18544.268
$ws.Cells.Item(22, 9).Value = 1351.3334
$ws.Cells.Item(22, 11).Value = 1351.3334
$ws.Cells.Item(22, 13).Value = -1056.3334
# Row 27
$ws.Cells.Item(27, 8).Value = 18544.268
$ws.Cells.Item(27, 9).Value = 1351.3334
$ws.Cells.Item(27, 11).Value = 1351.3334
$ws.Cells.Item(27, 13).Value = -1244.3334
# Row 46
$ws.Cells.Item(46, 8).Value = 2045.8334
$ws.Cells.Item(46, 10).Value = 3146
$ws.Cells.Item(46, 12).Value = 3146
$ws.Cells.Item(46, 14).Value = -3522
# Row 55
$ws.Cells.Item(55, 8).Value = 221
$ws.Cells.Item(55, 9).Value = 213.71428
$ws.Cells.Item(55, 11).Value = 213.71428
$ws.Cells.Item(55, 13).Value = -40.71428
# Row 68
$ws.Cells.Item(68, 8).Value = 3638.8333
$ws.Cells.Item(68, 9).Value = 3515.0908
$ws.Cells.Item(68, 10).Value = 5000
$ws.Cells.Item(68, 11).Value = 3515.0908
$ws.Cells.Item(68, 12).Value = 5000
$ws.Cells.Item(68, 13).Value = -2766.0908
$ws.Cells.Item(68, 14).Value = -6498
# Row 71
$ws.Cells.Item(71, 8).Value = 3638.8333
$ws.Cells.Item(71, 9).Value = 3515.0908
$ws.Cells.Item(71, 10).Value = 5000
$ws.Cells.Item(71, 11).Value = 17575.454
$ws.Cells.Item(71, 12).Value = 25000
$ws.Cells.Item(71, 13).Value = -13831.454
$ws.Cells.Item(71, 14).Value = -32488
# Row 122
$ws.Cells.Item(122, 8).Value = 21996.223
$ws.Cells.Item(122, 9).Value = 14125.842
$ws.Cells.Item(122, 11).Value = 42377.526
$ws.Cells.Item(122, 13).Value = -39927.526

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Cells.Item(54, 8).Value = 22371.223
$ws.Cells.Item(54, 10).Value = 24843.75
$ws.Cells.Item(54, 12).Value = 24843.75
$ws.Cells.Item(54, 14).Value = -25883.75
# Row 99
$ws.Cells.Item(99, 8).Value = 35499.5
$ws.Cells.Item(99, 9).Value = 26000
$ws.Cells.Item(99, 11).Value = 26000
$ws.Cells.Item(99, 13).Value = -23005
# Row 122
$ws.Cells.Item(122, 8).Value = 2001.2307
$ws.Cells.Item(122, 9).Value = 2002.5714
$ws.Cells.Item(122, 10).Value = 1995.6
$ws.Cells.Item(122, 11).Value = 6007.7142
$ws.Cells.Item(122, 12).Value = 5986.799999999999
$ws.Cells.Item(122, 13).Value = -3557.7142
$ws.Cells.Item(122, 14).Value = -10886.8
# Row 132
$ws.Cells.Item(132, 8).Value = 5952.4614
$ws.Cells.Item(132, 9).Value = 5952.4614
$ws.Cells.Item(132, 11).Value = 17857.3842
$ws.Cells.Item(132, 13).Value = -15327.3842
